# Update the "excel2skos" example worksheet: the concept URIs in the
# "Feuil2" sheet used a "#" fragment separator (e.g. .../days#monday);
# they should use a "/" path separator instead (.../days/monday).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil2")

# Column A (rows 8-15) holds the individual concept URIs.
$ws.Cells.Item(8, 1).Value = "http://data.sparna.fr/vocabularies/days/monday"
$ws.Cells.Item(9, 1).Value = "http://data.sparna.fr/vocabularies/days/tuesday"
$ws.Cells.Item(10, 1).Value = "http://data.sparna.fr/vocabularies/days/wednesday"
$ws.Cells.Item(11, 1).Value = "http://data.sparna.fr/vocabularies/days/thursday"
$ws.Cells.Item(12, 1).Value = "http://data.sparna.fr/vocabularies/days/friday"
$ws.Cells.Item(13, 1).Value = "http://data.sparna.fr/vocabularies/days/saturday"
$ws.Cells.Item(14, 1).Value = "http://data.sparna.fr/vocabularies/days/sunday"
$ws.Cells.Item(15, 1).Value = "http://data.sparna.fr/vocabularies/days/week-end"

# D15 references the Saturday/Sunday URIs for the "week-end" row's narrower list.
$ws.Cells.Item(15, 4).Value = "http://data.sparna.fr/vocabularies/days/saturday, http://data.sparna.fr/vocabularies/days/sunday"

# A10 (Wednesday) got touched in the fill-color dialog (explicitly set to "No
# Fill"), which stamps a new, otherwise-identical cell style on it.
$ws.Cells.Item(10, 1).Interior.ColorIndex = -4142

# Leave the selection on the last-edited cell, A10.
$ws.Activate() | Out-Null
$ws.Range("A10").Select() | Out-Null
